$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell while keeping its stored
# type as text (inline/shared string) instead of letting Excel infer a
# number, and keep the cell's style index back at the default (no style
# change), matching the workbook's original "plain General, no explicit
# number format" cells.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '87.747.80'
$ws.Range("E2").Value = '  -1.78%  '
Set-TextValue $ws.Range("D3") '3.099.07'
$ws.Range("E3").Value = '  -1.77%  '
$ws.Range("E4").Value = '  -0.39%  '
Set-TextValue $ws.Range("D5") '213.65'
$ws.Range("E5").Value = '  +1.89%  '
Set-TextValue $ws.Range("D6") '634.13'
$ws.Range("E6").Value = '  +4.01%  '
Set-TextValue $ws.Range("D7") '0.386'
$ws.Range("E7").Value = '  +0.98%  '
Set-TextValue $ws.Range("D8") '0.791'
$ws.Range("E8").Value = '  +16.41%  '
$ws.Range("E9").Value = '  -0.04%  '
Set-TextValue $ws.Range("D10") '3.092.98'
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("E11").Value = '  -1.27%  '
$ws.Range("E12").Value = '  +1.19%  '
$ws.Range("E13").Value = '  -0.30%  '
Set-TextValue $ws.Range("D14") '5.36'
$ws.Range("E14").Value = '  +3.66%  '
Set-TextValue $ws.Range("D15") '87.693.36'
$ws.Range("E15").Value = '  -1.98%  '
Set-TextValue $ws.Range("D16") '3.674.27'
$ws.Range("E16").Value = '  -1.90%  '
Set-TextValue $ws.Range("D17") '31.87'
$ws.Range("E17").Value = '  -0.98%  '
Set-TextValue $ws.Range("D18") '3.117.80'
$ws.Range("E18").Value = '  -2.13%  '
Set-TextValue $ws.Range("D19") '3.35'
$ws.Range("E19").Value = '  +4.44%  '
$ws.Range("E20").Value = '  +18.88%  '
$ws.Range("E21").Value = '  -0.57%  '
Set-TextValue $ws.Range("D22") '420.44'
$ws.Range("E22").Value = '  -2.86%  '
$ws.Range("E23").Value = '  -1.08%  '
Set-TextValue $ws.Range("D24") '4.86'
$ws.Range("E24").Value = '  -3.01%  '
$ws.Range("E25").Value = '  +7.23%  '
Set-TextValue $ws.Range("D26") '81.78'
$ws.Range("E26").Value = '  +10.46%  '
Set-TextValue $ws.Range("D27") '11.38'
$ws.Range("E27").Value = '  -1.20%  '
Set-TextValue $ws.Range("D29") '0.999'
$ws.Range("E29").Value = '  -0.03%  '
Set-TextValue $ws.Range("D30") '1.00'
$ws.Range("E30").Value = '  +0.29%  '
Set-TextValue $ws.Range("D31") '0.154'
$ws.Range("E31").Value = '  -7.72%  '
Set-TextValue $ws.Range("D32") '4.03'
$ws.Range("E32").Value = '  +1.48%  '
Set-TextValue $ws.Range("D33") '8.11'
$ws.Range("E33").Value = '  -2.15%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D34") '0.147'
$ws.Range("E34").Value = '  +15.88%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D35") '499.74'
$ws.Range("E35").Value = '  -4.42%  '
$ws.Range("E36").Value = '  -0.20%  '
Set-TextValue $ws.Range("D37") '1.27'
$ws.Range("E37").Value = '  +2.14%  '
$ws.Range("E38").Value = '  -0.88%  '
Set-TextValue $ws.Range("D39") '22.05'
$ws.Range("E39").Value = '  +1.88%  '
$ws.Range("E40").Value = '  -0.40%  '
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  -1.76%  '
Set-TextValue $ws.Range("D44") '1.82'
$ws.Range("E44").Value = '  -2.44%  '
$ws.Range("E45").Value = '  +9.66%  '
Set-TextValue $ws.Range("D46") '145.69'
$ws.Range("E46").Value = '  -2.30%  '
Set-TextValue $ws.Range("D47") '43.58'
$ws.Range("E47").Value = '  -1.51%  '
Set-TextValue $ws.Range("D48") '0.0647'
$ws.Range("E48").Value = '  +12.22%  '
Set-TextValue $ws.Range("D49") '161.12'
$ws.Range("E49").Value = '  -4.87%  '
Set-TextValue $ws.Range("D50") '0.710'
$ws.Range("E50").Value = '  +1.57%  '
$ws.Range("E51").Value = '  -3.29%  '
